# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet (copied from the existing "2022-Q3"
#    sheet so it keeps the same layout/formatting) right after "总计"
#  - fill it in with the new quarter's fund-holdings data
#  - insert a matching new row at the top of the "总计" summary sheet

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)        # "总计"
$q3 = $wb.Worksheets.Item(2)              # "2022-Q3" (template for the new sheet)

# --- create the new "2022-Q4" sheet, right after "总计" ---------------------
$q3.Copy($null, $summary)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4.Range("B2").Value = "'020023"
$q4.Range("C2").Value = "国泰事件驱动策略混合A"
$q4.Range("D2").Value = "'2.19"
$q4.Range("E2").Value = "'82.03"
$q4.Range("F2").Value = "'1.87"
$q4.Range("G2").Value = "'0.0410"
$q4.Range("H2").Value = 9

$q4.Range("B3").Value = "'015592"
$q4.Range("C3").Value = "国泰事件驱动策略混合C"
$q4.Range("D3").Value = "'0.01"
$q4.Range("E3").Value = "'82.03"
$q4.Range("F3").Value = "'1.87"
$q4.Range("G3").Value = "'0.0002"
$q4.Range("H3").Value = 9

# --- insert the new row into "总计", keeping its formatting ----------------
$summary.Rows(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.04

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# restore the originally-selected tab (last sheet, "2021-Q1")
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
